$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three oldest year rows (2004, 2008, 2009) -- rows 2,3,4.
# This shifts the remaining data rows (2010..2020) up to rows 2..12.
$ws.Range("A2:A4").EntireRow.Delete()

# Append the new 2021 data row at row 13, copying the label formatting
# (border/bold/alignment) used by the other year cells in column A.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("A13").Value = "2021年"
$rowVals = @(5069,7403,1651,14247,16373,5807,2598,4206,35082,6890,497,7787,7707,89627,3005,22010,3214,1257,32056,19635)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $col = [char]([int][char]'B' + $i)
    $ws.Range("$col`13").Value = $rowVals[$i]
}
